# Apply the edit described by the diff:
#  - Slide 2: swap the "Galleries" / "Testimonials" labels
#  - Slide 2: delete the extra "Reorder Icon" shape (id=77) at
#    off x=3435317, y=671073

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Swap the two tab labels.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Testimonials"
$s.Shapes.Item(4).TextFrame.TextRange.Text = "Galleries"

# Remove the stray "Reorder Icon" shape (id=77) that was left behind.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Reorder Icon" -and $shp.Id -eq 77) {
        $shp.Delete()
        break
    }
}
